# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the worker/period detail table (rows 16-25) on "Hoja1":
#   - Two new single-period workers (CRISTIAN DAVID RAMOS PEREZ / period 2205,
#     LEONARDO ANTONIO GARCIA MONROY / period 2205, with LEONARDO's document
#     number corrected to 20364364).
#   - The two four-period workers (DAVID RICARDO POLO PALENCIA and
#     DIEGO ANDRES HERNANDEZ LUNA) now cover periods 2301-2304 in ascending
#     order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$rows = @(
    @{ Row=16; TipoDoc="CC"; Doc="20364216";   Nombre="CRISTIAN DAVID RAMOS PEREZ";     Periodo="2205"; Mora=44000;  Salario=1100000 },
    @{ Row=17; TipoDoc="CC"; Doc="20364364";   Nombre="LEONARDO ANTONIO GARCIA MONROY"; Periodo="2205"; Mora=72000;  Salario=1800000 },
    @{ Row=18; TipoDoc="CC"; Doc="1102867201"; Nombre="DAVID RICARDO POLO PALENCIA";    Periodo="2301"; Mora=22000;  Salario=3300000 },
    @{ Row=19; TipoDoc="CC"; Doc="1100339566"; Nombre="DIEGO ANDRES HERNANDEZ LUNA";    Periodo="2301"; Mora=15120;  Salario=2268000 },
    @{ Row=20; TipoDoc="CC"; Doc="1102867201"; Nombre="DAVID RICARDO POLO PALENCIA";    Periodo="2302"; Mora=132000; Salario=3300000 },
    @{ Row=21; TipoDoc="CC"; Doc="1100339566"; Nombre="DIEGO ANDRES HERNANDEZ LUNA";    Periodo="2302"; Mora=90720;  Salario=2268000 },
    @{ Row=22; TipoDoc="CC"; Doc="1102867201"; Nombre="DAVID RICARDO POLO PALENCIA";    Periodo="2303"; Mora=132000; Salario=3300000 },
    @{ Row=23; TipoDoc="CC"; Doc="1100339566"; Nombre="DIEGO ANDRES HERNANDEZ LUNA";    Periodo="2303"; Mora=90720;  Salario=2268000 },
    @{ Row=24; TipoDoc="CC"; Doc="1102867201"; Nombre="DAVID RICARDO POLO PALENCIA";    Periodo="2304"; Mora=132000; Salario=3300000 },
    @{ Row=25; TipoDoc="CC"; Doc="1100339566"; Nombre="DIEGO ANDRES HERNANDEZ LUNA";    Periodo="2304"; Mora=90720;  Salario=2268000 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("B$i").Value = $r.TipoDoc
    $ws.Range("C$i").Value = $r.Doc
    $ws.Range("D$i").Value = $r.Nombre
    $ws.Range("E$i").Value = $r.Periodo
    $ws.Range("F$i").Value = $r.Mora
    $ws.Range("G$i").Value = $r.Salario
}
